$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "307.89"
Set-TextValue "E2" "-0.63%"
Set-TextValue "D3" "37.23"
Set-TextValue "E3" "-0.19%"
Set-TextValue "D4" "5.121"
Set-TextValue "E4" "-0.16%"
Set-TextValue "D5" "0.07818"
Set-TextValue "E5" "0.60%"
Set-TextValue "D6" "4.405"
Set-TextValue "E6" "-0.04%"
Set-TextValue "D7" "8.254"
Set-TextValue "E7" "0.78%"
Set-TextValue "D8" "1.877"
Set-TextValue "E8" "0.20%"
Set-TextValue "E9" "2.11%"
Set-TextValue "D10" "0.9215"
Set-TextValue "E10" "0.39%"
Set-TextValue "D11" "0.1079"
Set-TextValue "E11" "-9.42%"
Set-TextValue "D12" "0.1892"
Set-TextValue "E12" "-0.21%"
Set-TextValue "D13" "0.08883"
Set-TextValue "E13" "-5.61%"
Set-TextValue "D14" "0.03314"
Set-TextValue "E14" "-2.71%"
Set-TextValue "D15" "0.09589"
Set-TextValue "E15" "-0.88%"
Set-TextValue "D16" "0.001376"
Set-TextValue "E16" "-0.06%"
Set-TextValue "D17" "0.005709"
Set-TextValue "E17" "-1.07%"
Set-TextValue "D18" "3.396"
Set-TextValue "E18" "-3.94%"
Set-TextValue "E19" "0.70%"
Set-TextValue "D20" "6.312"
Set-TextValue "E20" "20.11%"
Set-TextValue "D21" "0.1287"
Set-TextValue "E21" "1.56%"
Set-TextValue "D22" "0.2416"
Set-TextValue "E22" "-6.72%"
Set-TextValue "D23" "0.04374"
Set-TextValue "E23" "1.14%"
Set-TextValue "D24" "0.001194"
Set-TextValue "E24" "-0.35%"
Set-TextValue "D25" "0.004269"
Set-TextValue "E25" "0.31%"
Set-TextValue "D26" "0.0001399"
Set-TextValue "E26" "7.59%"
Set-TextValue "E39" "4.56%"
Set-TextValue "D40" "0.05035"
Set-TextValue "E40" "0.26%"
Set-TextValue "D41" "0.007553"
Set-TextValue "E41" "-1.31%"
Set-TextValue "D42" "0.1351"
Set-TextValue "E42" "0.61%"
Set-TextValue "D43" "0.008639"
Set-TextValue "E43" "-12.03%"
Set-TextValue "D44" "0.002070"
Set-TextValue "E44" "-4.66%"
Set-TextValue "D45" "0.007896"
Set-TextValue "E45" "-9.56%"
Set-TextValue "D46" "0.00006511"
Set-TextValue "E46" "-2.92%"
Set-TextValue "E47" "-0.05%"
Set-TextValue "D48" "0.003295"
Set-TextValue "E48" "12.31%"
Set-TextValue "E49" "-16.48%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.05%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "-0.05%"
